$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits: "Start modeling ERD diagram" -> "Start modeling CDM diagram" ---
$ws.Range("B7").Value = "Start modeling CDM diagram"

# --- "Create Tables (...)" -> "Clean data and Create Tables (...)" ---
$ws.Range("B9").Value  = "Clean data and Create Tables ( Ridhwan Ibrahim) "
$ws.Range("B10").Value = "Clean data and Create Tables ( Giorgos Stefanis) "
$ws.Range("B11").Value = "Clean data and Create Tables (Michael Cipriani)"
$ws.Range("B12").Value = "Clean data and Create Tables( Youlun Wang)"
$ws.Range("B13").Value = "Clean data and Create Tables ( Seeram Govindan)"

# --- "procedures to load tables (...)" -> "stored procedures to load tables (...)" ---
$ws.Range("B14").Value = "stored procedures to load tables ( Hasan Husseini)"
$ws.Range("B15").Value = "stored procedures to load tables ( Ridhwan Ibrahim) "
$ws.Range("B16").Value = "stored procedures to load tables ( Giorgos Stefanis) "
$ws.Range("B17").Value = "stored procedures to load tables (Michael Cipriani)"
$ws.Range("B18").Value = "stored procedures to load tables( Youlun Wang)"
$ws.Range("B19").Value = "stored procedures to load tables ( Seeram Govindan)"

# --- Actual Duration (column F) tweaks, rows 8-18 ---
$ws.Range("F8").Value  = 5
$ws.Range("F9").Value  = 4
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 4
$ws.Range("F18").Value = 5

# --- Remove the duplicated "procedures to load tables ( Hasan Husseini)" row (old row 20) ---
# Deleting it shifts the remaining rows (old 21-24) up by one.
$ws.Rows(20).Delete()

# --- Update the "Load all production tables" row (now row 20) ---
$ws.Range("B20").Value = "Load all production tables (Ridhwan)"

# --- Bump Actual Start (column E) for the last three data rows ---
$ws.Range("E20").Value = 6
$ws.Range("E21").Value = 6
$ws.Range("E22").Value = 6

# --- Re-anchor the footer-row conditional format that doesn't auto-shift with row deletion ---
$footerCF = $ws.Range("B24:BO24").FormatConditions
if ($footerCF.Count -ge 1) {
    $footerCF.Item(1).ModifyAppliesToRange($ws.Range("B23:BO23"))
}

# --- Widen column B ---
$ws.Columns("B").ColumnWidth = 48.57

# --- Selection / view state ---
$ws.Range("F22").Select() | Out-Null
